$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TransactionSheet")

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "pain"
$ws.Range("C2").Value = 14
$ws.Range("D2").Value = "'829078"
$ws.Range("F2").Value = "'115666"
$ws.Range("H2").Value = 2000
$ws.Range("I2").Value = 0.2
$ws.Range("J2").Value = "2023-12-15 20:33:51"
$ws.Range("K2").Value = "CR"

# Row 3
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = "pain"
$ws.Range("C3").Value = 14
$ws.Range("D3").Value = "'829078"
$ws.Range("F3").Value = "'115666"
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 0.1
$ws.Range("J3").Value = "2023-12-15 20:35:51"
$ws.Range("K3").Value = "CR"

# Row 4
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = "pain"
$ws.Range("C4").Value = 14
$ws.Range("D4").Value = "'115666"
$ws.Range("F4").Value = "'829078"
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 0.05
$ws.Range("J4").Value = "2023-12-15 20:36:34"
$ws.Range("K4").Value = "DR"

# Row 5
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = "pain"
$ws.Range("C5").Value = 14
$ws.Range("D5").Value = "'115666"
$ws.Range("F5").Value = "'829078"
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 0.05
$ws.Range("J5").Value = "2023-12-15 20:36:36"
$ws.Range("K5").Value = "DR"

# Row 6
$ws.Range("A6").Value = 0
$ws.Range("B6").Value = "pain"
$ws.Range("C6").Value = 14
$ws.Range("D6").Value = "'829078"
$ws.Range("F6").Value = "'115666"
$ws.Range("H6").Value = 200
$ws.Range("I6").Value = 0.02
$ws.Range("J6").Value = "2023-12-15 20:50:27"
$ws.Range("K6").Value = "CR"

# Row 7
$ws.Range("A7").Value = 0
$ws.Range("B7").Value = "pain"
$ws.Range("C7").Value = 14
$ws.Range("D7").Value = "'829078"
$ws.Range("F7").Value = "'115666"
$ws.Range("H7").Value = 200
$ws.Range("I7").Value = 0.02
$ws.Range("J7").Value = "2023-12-17 21:55:50"
$ws.Range("K7").Value = "CR"

# Row 8
$ws.Range("A8").Value = 0
$ws.Range("B8").Value = "pain"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = "'829078"
$ws.Range("F8").Value = "'115666"
$ws.Range("H8").Value = 500
$ws.Range("I8").Value = 0.05
$ws.Range("J8").Value = "2023-12-17 22:06:51"
$ws.Range("K8").Value = "CR"

# Row 9
$ws.Range("A9").Value = 0
$ws.Range("B9").Value = "pain"
$ws.Range("C9").Value = 14
$ws.Range("D9").Value = "'829078"
$ws.Range("F9").Value = "'115666"
$ws.Range("H9").Value = 10000
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = "2023-12-17 22:27:04"
$ws.Range("K9").Value = "CR"

# Row 10
$ws.Range("A10").Value = 0
$ws.Range("B10").Value = "pain"
$ws.Range("C10").Value = 14
$ws.Range("D10").Value = "'829078"
$ws.Range("F10").Value = "'115666"
$ws.Range("H10").Value = 1000
$ws.Range("I10").Value = 0.1
$ws.Range("J10").Value = "2023-12-17 22:53:57"
$ws.Range("K10").Value = "CR"

# Row 11
$ws.Range("A11").Value = 0
$ws.Range("B11").Value = "pain"
$ws.Range("C11").Value = 14
$ws.Range("D11").Value = "'829078"
$ws.Range("F11").Value = "'441524"
$ws.Range("H11").Value = 500
$ws.Range("I11").Value = 0.05
$ws.Range("J11").Value = "2023-12-17 22:54:54"
$ws.Range("K11").Value = "CR"

# Row 12
$ws.Range("A12").Value = 0
$ws.Range("B12").Value = "pain"
$ws.Range("C12").Value = 14
$ws.Range("D12").Value = "'829078"
$ws.Range("F12").Value = "'115666"
$ws.Range("H12").Value = 500
$ws.Range("I12").Value = 0.05
$ws.Range("J12").Value = "2023-12-18 12:59:35"
$ws.Range("K12").Value = "CR"

# Row 13
$ws.Range("A13").Value = 0
$ws.Range("B13").Value = "pain"
$ws.Range("C13").Value = 14
$ws.Range("D13").Value = "'829078"
$ws.Range("F13").Value = "'441524"
$ws.Range("H13").Value = 10000
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = "2023-12-18 13:20:54"
$ws.Range("K13").Value = "DR"
